$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Automatic daily refresh of the MeteoCat summary: extraction timestamps
# (DATA_EXTRACCIO) move forward for every station row, and a handful of
# observed metrics (humidity, pressure, temperature) tick by their latest
# reading. Percentage-looking values (e.g. "85%") are entered with a
# leading apostrophe so Excel keeps them as literal text instead of
# auto-converting them to a numeric percentage, matching how the other
# text cells in this column are stored.
$ws.Range("E2").Value = '2026-02-07 19:47:40'
$ws.Range("E3").Value = '2026-02-07 19:47:43'
$ws.Range("H3").Value = "'85%"
$ws.Range("E4").Value = '2026-02-07 19:47:45'
$ws.Range("J4").Value = '1003.9 hPa'
$ws.Range("E5").Value = '2026-02-07 19:47:48'
$ws.Range("J5").Value = '1003.8 hPa'
$ws.Range("E6").Value = '2026-02-07 19:47:50'
$ws.Range("E7").Value = '2026-02-07 19:47:53'
$ws.Range("J7").Value = '1005.0 hPa'
$ws.Range("E8").Value = '2026-02-07 19:47:55'
$ws.Range("O8").Value = '8.7 °C'
$ws.Range("E9").Value = '2026-02-07 19:47:57'
$ws.Range("E10").Value = '2026-02-07 19:48:00'
$ws.Range("H10").Value = "'85%"
$ws.Range("O10").Value = '10.5 °C'
$ws.Range("E11").Value = '2026-02-07 19:48:02'
$ws.Range("N11").Value = '0.5 °C 19:23 TU'
$ws.Range("O11").Value = '3.2 °C'
$ws.Range("E12").Value = '2026-02-07 19:48:04'
$ws.Range("O12").Value = '12.2 °C'
$ws.Range("E13").Value = '2026-02-07 19:48:06'
$ws.Range("N13").Value = '6.0 °C 19:29 TU'
$ws.Range("O13").Value = '11.1 °C'
$ws.Range("E14").Value = '2026-02-07 19:48:09'
$ws.Range("E15").Value = '2026-02-07 19:48:11'
$ws.Range("J15").Value = '1004.1 hPa'
$ws.Range("E16").Value = '2026-02-07 19:48:14'
$ws.Range("E17").Value = '2026-02-07 19:48:16'
$ws.Range("E18").Value = '2026-02-07 19:48:19'
$ws.Range("H18").Value = "'94%"
$ws.Range("E19").Value = '2026-02-07 19:48:21'
$ws.Range("H19").Value = "'78%"
$ws.Range("O19").Value = '7.5 °C'
$ws.Range("E20").Value = '2026-02-07 19:48:24'
$ws.Range("E21").Value = '2026-02-07 19:48:26'
$ws.Range("J21").Value = '1004.0 hPa'
$ws.Range("O21").Value = '8.7 °C'
$ws.Range("E22").Value = '2026-02-07 19:48:29'
$ws.Range("E23").Value = '2026-02-07 19:48:31'
$ws.Range("J23").Value = '1003.8 hPa'
$ws.Range("O23").Value = '10.2 °C'
$ws.Range("E24").Value = '2026-02-07 19:48:34'
$ws.Range("E25").Value = '2026-02-07 19:48:36'
$ws.Range("H25").Value = "'83%"
$ws.Range("E26").Value = '2026-02-07 19:48:38'
$ws.Range("E27").Value = '2026-02-07 19:48:41'
$ws.Range("E28").Value = '2026-02-07 19:48:43'
$ws.Range("J28").Value = '1005.8 hPa'
$ws.Range("O28").Value = '4.9 °C'
$ws.Range("E29").Value = '2026-02-07 19:48:45'
$ws.Range("E30").Value = '2026-02-07 19:48:48'
$ws.Range("H30").Value = "'68%"
$ws.Range("O30").Value = '-4.4 °C'
$ws.Range("E31").Value = '2026-02-07 19:48:50'
$ws.Range("E32").Value = '2026-02-07 19:48:53'
$ws.Range("H32").Value = "'54%"
$ws.Range("E33").Value = '2026-02-07 19:48:55'
$ws.Range("E34").Value = '2026-02-07 19:48:58'
$ws.Range("O34").Value = '7.8 °C'
$ws.Range("E35").Value = '2026-02-07 19:49:00'
$ws.Range("E36").Value = '2026-02-07 19:49:03'
